# feat: add edge 0-shot
# Adds two summary rows below the existing data/average row:
#   row 106: STDEV.S() of each metric column (B:F)
#   row 107: the 95% CI half-width (StdErr * 1.96) for each metric column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 106 - sample standard deviation of each column (B..F)
$ws.Range("B106").Formula = "=STDEV.S(B2:B104)"
$ws.Range("C106:F106").Formula = "=STDEV.S(C2:C104)"

# Row 107 - 95% confidence interval half width based on row 106
$ws.Range("B107").Formula = "=B106/SQRT(103)*1.96"
$ws.Range("C107:F107").Formula = "=C106/SQRT(103)*1.96"

# Move the selection to match the post-edit view (top of sheet, column F)
$ws.Range("F1").Select() | Out-Null
